# Commit: "Added Ford Model T"
# Appends a new data row (row 53) for the "Ford Model T" vehicle to Sheet1,
# mirroring the layout/formulas used by the existing rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 53
$prevRow = $row - 1

# --- Plain data values -------------------------------------------------
$ws.Cells.Item($row, 1).Value = "Ford Model T"   # A53 - Vehicle
$ws.Cells.Item($row, 2).Value = 1908             # B53 - Intro Year
$ws.Cells.Item($row, 3).Value = 1                # C53 - Year Order
$ws.Cells.Item($row, 4).Value = "Light Goods"    # D53 - Vehicle Type
$ws.Cells.Item($row, 6).Value = 40               # F53 - Top Speed
$ws.Cells.Item($row, 7).Value = 4                # G53 - Capacity (goods)
$ws.Cells.Item($row, 10).Value = "x"             # J53 - Done

# --- Formulas (same pattern as the row above) ---------------------------
$ws.Cells.Item($row, 5).Formula = "=IF(B$row > 1900, ((B$row-1900)*10)+400+C$row, ((B$row-1730)*2)+C$row)+VLOOKUP(D$row,'ID Scheme'!`$A`$2:`$B`$6,2, FALSE)"
$ws.Cells.Item($row, 8).Formula = "=SQRT(F$row*G$row)/`$B`$1"
$ws.Cells.Item($row, 9).Formula = "=H$row*0.9"

# --- Carry over the number formatting used by Cost / Running Cost / Done
$ws.Cells.Item($row, 8).NumberFormat = $ws.Cells.Item($prevRow, 8).NumberFormat
$ws.Cells.Item($row, 9).NumberFormat = $ws.Cells.Item($prevRow, 9).NumberFormat
$ws.Cells.Item($row, 10).NumberFormat = $ws.Cells.Item($prevRow, 10).NumberFormat

# --- Update the view so the new row is visible/selected, like the author left it
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("I53").Select()

$wb.Save()
